$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

# --- Row 51 ---
$ws.Range("I51").NumberFormat = "@"
$ws.Range("Y51").NumberFormat = "@"
$ws.Range("AA51").NumberFormat = "@"
$ws.Range("A51").Value = 131158283
$ws.Range("B51").Value = 58043
$ws.Range("D51").Value = "NT"
$ws.Range("E51").Value = 103021
$ws.Range("F51").Value = "Talltita"
$ws.Range("G51").Value = "Poecile montanus"
$ws.Range("H51").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I51").Value = "1"
$ws.Range("M51").Value = "upprörd, varnande"
$ws.Range("P51").Value = "Hundviken, Dlr"
$ws.Range("Q51").Value = 525550
$ws.Range("R51").Value = 6716117
$ws.Range("S51").Value = 5
$ws.Range("T51").Value = "Dalarna"
$ws.Range("U51").Value = "Falun"
$ws.Range("V51").Value = "Dalarna"
$ws.Range("W51").Value = "Aspeboda"
$ws.Range("Y51").Value = "2026-02-01"
$ws.Range("AA51").Value = "2026-02-01"
$ws.Range("AC51").Value = "Hört fågeln och registrerat med app."
$ws.Range("AD51").Value = $false
$ws.Range("AE51").Value = $false
$ws.Range("AG51").Value = $false
$ws.Range("AW51").Value = "Daniel Alexandersson"
$ws.Range("AX51").Value = "Daniel Alexandersson"

# --- Row 52 ---
$ws.Range("I52").NumberFormat = "@"
$ws.Range("Y52").NumberFormat = "@"
$ws.Range("AA52").NumberFormat = "@"
$ws.Range("A52").Value = 131158205
$ws.Range("B52").Value = 58043
$ws.Range("D52").Value = "NT"
$ws.Range("E52").Value = 103021
$ws.Range("F52").Value = "Talltita"
$ws.Range("G52").Value = "Poecile montanus"
$ws.Range("H52").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I52").Value = "1"
$ws.Range("M52").Value = "upprörd, varnande"
$ws.Range("P52").Value = "Hundviken, Dlr"
$ws.Range("Q52").Value = 525731
$ws.Range("R52").Value = 6716197
$ws.Range("S52").Value = 5
$ws.Range("T52").Value = "Dalarna"
$ws.Range("U52").Value = "Falun"
$ws.Range("V52").Value = "Dalarna"
$ws.Range("W52").Value = "Aspeboda"
$ws.Range("Y52").Value = "2026-02-08"
$ws.Range("AA52").Value = "2026-02-08"
$ws.Range("AD52").Value = $false
$ws.Range("AE52").Value = $false
$ws.Range("AG52").Value = $false
$ws.Range("AW52").Value = "Daniel Alexandersson"
$ws.Range("AX52").Value = "Daniel Alexandersson"

# --- Row 53 ---
$ws.Range("I53").NumberFormat = "@"
$ws.Range("Y53").NumberFormat = "@"
$ws.Range("AA53").NumberFormat = "@"
$ws.Range("A53").Value = 131158127
$ws.Range("B53").Value = 58043
$ws.Range("D53").Value = "NT"
$ws.Range("E53").Value = 103021
$ws.Range("F53").Value = "Talltita"
$ws.Range("G53").Value = "Poecile montanus"
$ws.Range("H53").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I53").Value = "1"
$ws.Range("M53").Value = "upprörd, varnande"
$ws.Range("P53").Value = "Hundviken, Dlr"
$ws.Range("Q53").Value = 525695
$ws.Range("R53").Value = 6716177
$ws.Range("S53").Value = 5
$ws.Range("T53").Value = "Dalarna"
$ws.Range("U53").Value = "Falun"
$ws.Range("V53").Value = "Dalarna"
$ws.Range("W53").Value = "Aspeboda"
$ws.Range("Y53").Value = "2026-01-11"
$ws.Range("AA53").Value = "2026-01-11"
$ws.Range("AD53").Value = $false
$ws.Range("AE53").Value = $false
$ws.Range("AG53").Value = $false
$ws.Range("AW53").Value = "Daniel Alexandersson"
$ws.Range("AX53").Value = "Daniel Alexandersson"

# --- Row 54 ---
$ws.Range("I54").NumberFormat = "@"
$ws.Range("Y54").NumberFormat = "@"
$ws.Range("AA54").NumberFormat = "@"
$ws.Range("A54").Value = 131158214
$ws.Range("B54").Value = 58043
$ws.Range("D54").Value = "NT"
$ws.Range("E54").Value = 103021
$ws.Range("F54").Value = "Talltita"
$ws.Range("G54").Value = "Poecile montanus"
$ws.Range("H54").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("I54").Value = "1"
$ws.Range("M54").Value = "upprörd, varnande"
$ws.Range("P54").Value = "Hundviken, Dlr"
$ws.Range("Q54").Value = 525566
$ws.Range("R54").Value = 6716104
$ws.Range("S54").Value = 5
$ws.Range("T54").Value = "Dalarna"
$ws.Range("U54").Value = "Falun"
$ws.Range("V54").Value = "Dalarna"
$ws.Range("W54").Value = "Aspeboda"
$ws.Range("Y54").Value = "2026-01-25"
$ws.Range("AA54").Value = "2026-01-25"
$ws.Range("AD54").Value = $false
$ws.Range("AE54").Value = $false
$ws.Range("AG54").Value = $false
$ws.Range("AW54").Value = "Daniel Alexandersson"
$ws.Range("AX54").Value = "Daniel Alexandersson"

Write-Output "done"
